$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New Pac-12 rows (55-66). Shared-string order matters for a faithful replay:
# author filled column A (school) top-to-bottom, then column B (conference),
# then column C (filename) mostly alphabetically with the "*_state.py" files
# appended afterwards, then column H (football_url) in row order.
# ---------------------------------------------------------------------------

# -- Column A : school --------------------------------------------------
$ws.Range("A55").Value = "Arizona"
$ws.Range("A56").Value = "Arizona State"
$ws.Range("A57").Value = "Cal"
$ws.Range("A58").Value = "Colorado"
$ws.Range("A59").Value = "Oregon"
$ws.Range("A60").Value = "Oregon State"
$ws.Range("A61").Value = "Stanford"
$ws.Range("A62").Value = "UCLA"
$ws.Range("A63").Value = "USC"
$ws.Range("A64").Value = "Utah"
$ws.Range("A65").Value = "Washington"
$ws.Range("A66").Value = "Washington State"

# -- Column B : conference ----------------------------------------------
$ws.Range("B55:B66").Value = "Pac-12"

# -- Column C : filename --------------------------------------------------
$ws.Range("C55").Value = "arizona.py"
$ws.Range("C57").Value = "cal.py"
$ws.Range("C58").Value = "colorado.py"
$ws.Range("C59").Value = "oregon.py"
$ws.Range("C61").Value = "stanford.py"
$ws.Range("C62").Value = "ucla.py"
$ws.Range("C63").Value = "usc.py"
$ws.Range("C64").Value = "utah.py"
$ws.Range("C65").Value = "washington.py"
$ws.Range("C56").Value = "arizona_state.py"
$ws.Range("C60").Value = "oregon_state.py"
$ws.Range("C66").Value = "washington_state.py"

# -- Columns D/E/F : grid / table / ul flags (all new rows are "ul") -----
$ws.Range("D55:D66").Value = 0
$ws.Range("E55:E66").Value = 1
$ws.Range("F55:F66").Value = 0

# -- Column G : scrape_date (reuse the existing date style/format) -------
$ws.Range("G54").Copy()
$ws.Range("G55:G66").PasteSpecial(-4122)
$ws.Range("G55").Value = (Get-Date -Year 2018 -Month 4 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G56").Value = (Get-Date -Year 2018 -Month 4 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G57").Value = (Get-Date -Year 2018 -Month 5 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G58").Value = (Get-Date -Year 2018 -Month 5 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G59").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G60").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G61").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G62").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G63").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G64").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G65").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G66").Value = (Get-Date -Year 2018 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0)

# -- Column H : football_url (text + hyperlink), in row order ------------
$ws.Range("H55").Value = "http://arizonawildcats.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H55"), "http://arizonawildcats.com/roster.aspx?path=football")

$ws.Range("H56").Value = "http://thesundevils.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H56"), "http://thesundevils.com/roster.aspx?path=football")

$ws.Range("H57").Value = "http://calbears.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H57"), "http://calbears.com/roster.aspx?path=football")

$ws.Range("H58").Value = "http://cubuffs.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H58"), "http://cubuffs.com/roster.aspx?path=football")

$ws.Range("H59").Value = "http://goducks.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H59"), "http://goducks.com/roster.aspx?path=football")

$ws.Range("H60").Value = "http://osubeavers.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H60"), "http://osubeavers.com/roster.aspx?path=football")

$ws.Range("H61").Value = "http://gostanford.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H61"), "http://gostanford.com/roster.aspx?path=football")

$ws.Range("H62").Value = "http://uclabruins.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H62"), "http://uclabruins.com/roster.aspx?path=football")

$ws.Range("H63").Value = "http://usctrojans.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H63"), "http://usctrojans.com/roster.aspx?path=football")

$ws.Range("H64").Value = "http://utahutes.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H64"), "http://utahutes.com/roster.aspx?path=football")

$ws.Range("H65").Value = "http://gohuskies.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H65"), "http://gohuskies.com/roster.aspx?path=football")

$ws.Range("H66").Value = "http://wsucougars.com/roster.aspx?path=football "
$ws.Hyperlinks.Add($ws.Range("H66"), "http://wsucougars.com/roster.aspx?path=football")

# ---------------------------------------------------------------------------
# Column widths: A and C grew slightly to fit the new, longer Pac-12 values.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 15.833333333333334
$ws.Columns("C").ColumnWidth = 21.666666666666668

# ---------------------------------------------------------------------------
# View: scroll the frozen (bottom) pane down and move the selection to the
# new last row, mirroring where the author ended up after pasting the data.
# ---------------------------------------------------------------------------
$ws.Range("A67").Select()
